# Refresh ligand-receptor summary stats (Efna5-Ephb6) with updated TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2708946666666667
$ws.Range("H2").Value = 0.812684
$ws.Range("I2").Value = 0.1616296696421007
$ws.Range("J2").Value = 0.1616296696421007
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08324533333333334
$ws.Range("N2").Value = 0.249736
$ws.Range("O2").Value = 0.05078606388889115
$ws.Range("P2").Value = 0.05078606388889115
$ws.Range("Q2").Value = 0.02255071682488889
$ws.Range("R2").Value = 0.202956451424
$ws.Range("S2").Value = 0.008208534728784094
$ws.Range("T2").Value = 0.008208534728784094

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2708946666666667
$ws.Range("H3").Value = 0.812684
$ws.Range("I3").Value = 0.1616296696421007
$ws.Range("J3").Value = 0.1616296696421007
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4583163333333333
$ws.Range("N3").Value = 1.374949
$ws.Range("O3").Value = 0.2796082573516313
$ws.Range("P3").Value = 0.2796082573516313
$ws.Range("Q3").Value = 0.1241554503462222
$ws.Range("R3").Value = 1.117399053116
$ws.Range("S3").Value = 0.04519299026494763
$ws.Range("T3").Value = 0.04519299026494762

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2708946666666667
$ws.Range("H4").Value = 0.812684
$ws.Range("I4").Value = 0.1616296696421007
$ws.Range("J4").Value = 0.1616296696421007
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.097575666666667
$ws.Range("N4").Value = 3.292727
$ws.Range("O4").Value = 0.6696056787594775
$ws.Range("P4").Value = 0.6696056787594775
$ws.Range("Q4").Value = 0.2973273943631111
$ws.Range("R4").Value = 2.675946549268
$ws.Range("S4").Value = 0.1082281446483689
$ws.Range("T4").Value = 0.1082281446483689

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.058121333333333
$ws.Range("H5").Value = 3.174364
$ws.Range("I5").Value = 0.6313295261673385
$ws.Range("J5").Value = 0.6313295261673384
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08324533333333334
$ws.Range("N5").Value = 0.249736
$ws.Range("O5").Value = 0.05078606388889115
$ws.Range("P5").Value = 0.05078606388889115
$ws.Range("Q5").Value = 0.08808366310044446
$ws.Range("R5").Value = 0.7927529679040001
$ws.Range("S5").Value = 0.03206274165087783
$ws.Range("T5").Value = 0.03206274165087782

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.058121333333333
$ws.Range("H6").Value = 3.174364
$ws.Range("I6").Value = 0.6313295261673385
$ws.Range("J6").Value = 0.6313295261673384
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4583163333333333
$ws.Range("N6").Value = 1.374949
$ws.Range("O6").Value = 0.2796082573516313
$ws.Range("P6").Value = 0.2796082573516313
$ws.Range("Q6").Value = 0.4849542897151112
$ws.Range("R6").Value = 4.364588607436001
$ws.Range("S6").Value = 0.1765249486262806
$ws.Range("T6").Value = 0.1765249486262806

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.058121333333333
$ws.Range("H7").Value = 3.174364
$ws.Range("I7").Value = 0.6313295261673385
$ws.Range("J7").Value = 0.6313295261673384
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.097575666666667
$ws.Range("N7").Value = 3.292727
$ws.Range("O7").Value = 0.6696056787594775
$ws.Range("P7").Value = 0.6696056787594775
$ws.Range("Q7").Value = 1.161368227847556
$ws.Range("R7").Value = 10.452314050628
$ws.Range("S7").Value = 0.4227418358901801
$ws.Range("T7").Value = 0.4227418358901799

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3470046666666667
$ws.Range("H8").Value = 1.041014
$ws.Range("I8").Value = 0.2070408041905609
$ws.Range("J8").Value = 0.2070408041905609
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08324533333333334
$ws.Range("N8").Value = 0.249736
$ws.Range("O8").Value = 0.05078606388889115
$ws.Range("P8").Value = 0.05078606388889115
$ws.Range("Q8").Value = 0.02888651914488889
$ws.Range("R8").Value = 0.259978672304
$ws.Range("S8").Value = 0.01051478750922923
$ws.Range("T8").Value = 0.01051478750922923

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3470046666666667
$ws.Range("H9").Value = 1.041014
$ws.Range("I9").Value = 0.2070408041905609
$ws.Range("J9").Value = 0.2070408041905609
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4583163333333333
$ws.Range("N9").Value = 1.374949
$ws.Range("O9").Value = 0.2796082573516313
$ws.Range("P9").Value = 0.2796082573516313
$ws.Range("Q9").Value = 0.1590379064762222
$ws.Range("R9").Value = 1.431341158286
$ws.Range("S9").Value = 0.05789031846040306
$ws.Range("T9").Value = 0.05789031846040305

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.3470046666666667
$ws.Range("H10").Value = 1.041014
$ws.Range("I10").Value = 0.2070408041905609
$ws.Range("J10").Value = 0.2070408041905609
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.097575666666667
$ws.Range("N10").Value = 3.292727
$ws.Range("O10").Value = 0.6696056787594775
$ws.Range("P10").Value = 0.6696056787594775
$ws.Range("Q10").Value = 0.3808638783531111
$ws.Range("R10").Value = 3.427774905178
$ws.Range("S10").Value = 0.1386356982209286
$ws.Range("T10").Value = 0.1386356982209286
